$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the six existing fold-change values before we overwrite anything
$v1 = $ws.Range("B2").Value2
$v2 = $ws.Range("C2").Value2
$v3 = $ws.Range("D2").Value2
$v4 = $ws.Range("E2").Value2
$v5 = $ws.Range("F2").Value2
$v6 = $ws.Range("G2").Value2

# Drop the now-unused columns D:G (this also removes D1:G1 headers and D2:G2 values)
$ws.Range("D1:G1").EntireColumn.Delete()

# New header row
$ws.Range("B1").Value = "Condition"
$ws.Range("C1").Value = "Fold Change"

# Give rows 3-7 column A the same style (border/bold/alignment) as A2 by copying formats
$ws.Range("A2").Copy()
$ws.Range("A3:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row index values in column A
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

# Condition labels in column B
$ws.Range("B2").Value = "dusp11 -dox foldchange"
$ws.Range("B3").Value = "ifnb -dox foldchange"
$ws.Range("B4").Value = "mx1 -dox foldchange"
$ws.Range("B5").Value = "dusp11 +dox foldchange"
$ws.Range("B6").Value = "ifnb +dox foldchange"
$ws.Range("B7").Value = "mx1 +dox foldchange"

# Fold change values in column C
$ws.Range("C2").Value = $v1
$ws.Range("C3").Value = $v2
$ws.Range("C4").Value = $v3
$ws.Range("C5").Value = $v4
$ws.Range("C6").Value = $v5
$ws.Range("C7").Value = $v6
